$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Crit and Crit Fail Sounds" (row 20) shipped in 1.5.8:
#  - fill in the Completed Version column
#  - tweak the description wording
#  - hide the now-completed row (matches the other completed rows)
$ws.Range("C20").Value = "1.5.8"
$ws.Range("B20").Value = "When you roll minimum or maximum on a roll, do the willhelm scream or the air horn sound effects"
$ws.Rows(20).Hidden = $true

# Grow the autofilter down one row to include the new last data row,
# keeping the existing "blanks" filter on column C (Completed Version).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:E21").AutoFilter(3, @(""), 7)

# Keep the _FilterDatabase defined name in sync with the grown autofilter range.
$fd = $wb.Names.Item("Sheet1!_FilterDatabase")
$fd.RefersTo = "=Sheet1!`$A`$1:`$E`$21"

# Reflect where the user's cursor ended up after making the edits.
$null = $ws.Range("B26").Select()
